$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("08-10-2021", -1.57, -0.55, -0.21),
    @("12-10-2021", -1.92, -0.8, -0.35),
    @("13-10-2021", -2, -1.04, -0.41),
    @("14-10-2021", -2.06, -1.23, -0.48)
)

$startRow = 195
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    # Column A holds dates formatted as plain text ("DD-MM-YYYY"). Writing
    # them directly via .Value lets Excel's smart-parsing turn unambiguous
    # ones (day <= 12) into real dates. Route through a text formula first,
    # then convert the formula to its literal value in place (Copy +
    # PasteSpecial values-only) so the result is a plain shared-string cell
    # with the default style, matching how the rest of the column is stored.
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.Formula = '="' + $data[$i][0] + '"'
    $cellA.Copy()
    $cellA.PasteSpecial(-4163)

    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
